# Weekly update: insert 3 new rows (a new "Doctor Davis" batch dated 44617)
# above the existing data block that starts at row 162, pushing the rest of
# the sheet down by 3 rows (old 162:185 -> new 165:188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 162, shifting rows 162:185 down to 165:188.
$ws.Rows("162:164").Insert()

# Common values shared by all three new rows.
$mercadoId = 11
$mercado   = 'Vega Monumental Concepción'
$region    = 'Bíobío'
$fecha     = 44617
$codreg    = 8
$tipo      = 'Fruta'
$prodId    = 100103
$producto  = 'Frutos de hueso (carozo)'
$catId     = 100103004
$categoria = 'Durazno'
$variedad  = 'Doctor Davis'
$unidad    = '$/caja 15 kilos empedrada'
$origen    = "Región de O'Higgins"
$kgUnidad  = 15

# Row 162: Especial
$ws.Cells.Item(162, 1).Value  = $mercadoId
$ws.Cells.Item(162, 2).Value  = $mercado
$ws.Cells.Item(162, 3).Value  = $region
$ws.Cells.Item(162, 4).Value  = $fecha
$ws.Cells.Item(162, 5).Value  = $codreg
$ws.Cells.Item(162, 6).Value  = $tipo
$ws.Cells.Item(162, 7).Value  = $prodId
$ws.Cells.Item(162, 8).Value  = $producto
$ws.Cells.Item(162, 9).Value  = $catId
$ws.Cells.Item(162, 10).Value = $categoria
$ws.Cells.Item(162, 11).Value = $variedad
$ws.Cells.Item(162, 12).Value = "Especial"
$ws.Cells.Item(162, 13).Value = 50
$ws.Cells.Item(162, 14).Value = 13000
$ws.Cells.Item(162, 15).Value = 13000
$ws.Cells.Item(162, 16).Value = 13000
$ws.Cells.Item(162, 17).Value = $unidad
$ws.Cells.Item(162, 18).Value = $origen
$ws.Cells.Item(162, 19).Value = 867
$ws.Cells.Item(162, 20).Value = $kgUnidad

# Row 163: Primera
$ws.Cells.Item(163, 1).Value  = $mercadoId
$ws.Cells.Item(163, 2).Value  = $mercado
$ws.Cells.Item(163, 3).Value  = $region
$ws.Cells.Item(163, 4).Value  = $fecha
$ws.Cells.Item(163, 5).Value  = $codreg
$ws.Cells.Item(163, 6).Value  = $tipo
$ws.Cells.Item(163, 7).Value  = $prodId
$ws.Cells.Item(163, 8).Value  = $producto
$ws.Cells.Item(163, 9).Value  = $catId
$ws.Cells.Item(163, 10).Value = $categoria
$ws.Cells.Item(163, 11).Value = $variedad
$ws.Cells.Item(163, 12).Value = "Primera"
$ws.Cells.Item(163, 13).Value = 100
$ws.Cells.Item(163, 14).Value = 11000
$ws.Cells.Item(163, 15).Value = 11000
$ws.Cells.Item(163, 16).Value = 11000
$ws.Cells.Item(163, 17).Value = $unidad
$ws.Cells.Item(163, 18).Value = $origen
$ws.Cells.Item(163, 19).Value = 733
$ws.Cells.Item(163, 20).Value = $kgUnidad

# Row 164: Segunda
$ws.Cells.Item(164, 1).Value  = $mercadoId
$ws.Cells.Item(164, 2).Value  = $mercado
$ws.Cells.Item(164, 3).Value  = $region
$ws.Cells.Item(164, 4).Value  = $fecha
$ws.Cells.Item(164, 5).Value  = $codreg
$ws.Cells.Item(164, 6).Value  = $tipo
$ws.Cells.Item(164, 7).Value  = $prodId
$ws.Cells.Item(164, 8).Value  = $producto
$ws.Cells.Item(164, 9).Value  = $catId
$ws.Cells.Item(164, 10).Value = $categoria
$ws.Cells.Item(164, 11).Value = $variedad
$ws.Cells.Item(164, 12).Value = "Segunda"
$ws.Cells.Item(164, 13).Value = 100
$ws.Cells.Item(164, 14).Value = 9000
$ws.Cells.Item(164, 15).Value = 9000
$ws.Cells.Item(164, 16).Value = 9000
$ws.Cells.Item(164, 17).Value = $unidad
$ws.Cells.Item(164, 18).Value = $origen
$ws.Cells.Item(164, 19).Value = 600
$ws.Cells.Item(164, 20).Value = $kgUnidad
